# Adds a new "version" column (I) to the grade list: a bold "version"
# header in I2 and a repeating A/B/C/D value in I3:I18 (one of four grading
# "versions" per student row). Also nudges a couple of print/view settings
# to match the author's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "version" column ---------------------------------------------
$ws.Range("I2").Value = "version"
$ws.Range("I2").Font.Bold = $true

$letters = @("A", "B", "C", "D")
for ($r = 3; $r -le 18; $r++) {
    $ws.Range("I$r").Value = $letters[($r - 3) % 4]
}

# --- Print setup (paper size / orientation) ----------------------------
$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait

# --- Leave the cursor where the author left it after editing the list --
[void]$ws.Range("I22").Select()
